# Update "paises" (countries) stats and re-sort by total cases, per commit
# "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp banner.
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 13:15"

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# 2) Apply the refreshed per-country statistics (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
#    Rows below refer to the ORIGINAL (pre-sort) layout; the table gets
#    re-sorted by "Casos totales" afterwards, which naturally reshuffles any
#    rows whose updated totals change their rank (matches how the published
#    sheet re-ranks Nepal/Panama, Austria/Nigeria, Libia/Ghana,
#    Bosnia/Birmania, Islandia/Botsuana/Somalia, Gibraltar/Taiwan and
#    Liechtenstein/Bermudas).

Set-Row 4   8219123 2808 5320386 2675983 0 37  222754   # Estados Unidos
Set-Row 16  522387  4552 420910  71607   0 265 29870    # Iran
Set-Row 32  172516  4026 127076  39691   0 75  5749     # Rumania
Set-Row 40  126137  4392 88040   37382   0 21  715      # Nepal
Set-Row 43  112849  1412 104943  7451    0 3   455      # Emiratos Arabes Unidos
Set-Row 58  74422   3105 50500   21802   0 5   2120     # Suiza
Set-Row 63  61387   1163 47618   12887   0 5   882      # Austria
Set-Row 70  47845   1169 26062   21084   0 18  699      # Libia
Set-Row 76  40073   47   33516   5072    0 4   1485     # Afganistan
Set-Row 82  32845   621  24603   7262    0 8   980      # Bosnia y Herzegovina
Set-Row 93  18758   629  12259   6323    0 6   176      # Malasia
Set-Row 98  15368   20   13704   1347    0 1   317      # Senegal
Set-Row 139 4282    122  3142    1095    0 0   45       # Malta
Set-Row 144 3929    92   2713    1206    0 0   10       # Islandia
Set-Row 177 544     13   441     103     0 0   0        # Gibraltar
Set-Row 182 482     2    471     11      0 0   0        # Islas Feroe
Set-Row 193 192     9    132     59      0 0   1        # Liechtenstein

# 3) Re-sort the whole country table (A4:H220) by "Casos totales" (col B)
#    descending, same as the published sheet does after each refresh.
$sortRange = $ws.Range("A4:H220")
$sortKey = $ws.Range("B4:B220")
$sortRange.Sort($sortKey, 2)
